$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 and row 3 labels (A2, A3): "Mitglied..." -> "Gemeldete..."
$ws.Range("A2").Value = "Gemeldete Krankenkasse (vollständiger Name)"
$ws.Range("A3").Value = "Gemeldete Krankenkasse (Abkürzung)"

# Update B7 (Eintragungsdatum value) from "15.12.2023" to "01.01.2024" as text
$ws.Range("B7").Value = "01.01.2024"

# Update selection to B9 as in diff
$ws.Range("B9").Select()
